$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Replace the placeholder folder name 'day2-files' with the real
#    folder path 'FrontEnd/Code/Day 2' inside item 2 of the Day 2 steps.
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "day2-files", $true, $false, $false, $false, $false,
    $true, 1, $false, "FrontEnd/Code/Day 2", 2
)

# --------------------------------------------------------------------
# 2) Re-locate the freshly inserted text and drop the (single, auto-
#    managed) "_GoBack" bookmark right after the "D" of "Day" - this
#    mirrors where Word leaves it after the last text edit. Adding a
#    bookmark with this reserved name automatically relocates it from
#    wherever it previously lived (Word only ever keeps one).
# --------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute(
    "FrontEnd/Code/D", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$posBeforeFrontEnd = $find.Start
$posAfterD = $find.End

$goBackRange = $d.Range($posAfterD, $posAfterD)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --------------------------------------------------------------------
# 3) Break the remaining text into the same run boundaries as the
#    final document (the leading quote | "FrontEnd/Code/D" | "ay" |
#    " 2" | the trailing text) by placing and immediately removing
#    throw-away bookmarks at those boundaries.
# --------------------------------------------------------------------
$tmp0Range = $d.Range($posBeforeFrontEnd, $posBeforeFrontEnd)
$d.Bookmarks.Add("ZZZTempSplit0", $tmp0Range)

$posAfterAy = $posAfterD + 2
$tmp1Range = $d.Range($posAfterAy, $posAfterAy)
$d.Bookmarks.Add("ZZZTempSplit1", $tmp1Range)

$posAfterSpace2 = $posAfterAy + 2
$tmp2Range = $d.Range($posAfterSpace2, $posAfterSpace2)
$d.Bookmarks.Add("ZZZTempSplit2", $tmp2Range)

$d.Bookmarks("ZZZTempSplit0").Delete()
$d.Bookmarks("ZZZTempSplit1").Delete()
$d.Bookmarks("ZZZTempSplit2").Delete()
